# Update "想去人数" (interest counts) figures for two 漫展 (events)
# in both the "展览" sheet and the "全部类型" sheet, matching data
# refreshed at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 9641
    $ws.Range("F5").Value = 543
}
